$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Empty" labels first (matches the original authoring order so the
# shared-strings table comes out in the same sequence as the real edit).
$ws.Range("F1").Value = "Empty"
$ws.Range("F2").Value = "Empty{.empty}"

# Translate the header row from Chinese to English.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Complex"
$ws.Range("E1").Value = "Ignored"

# Translate the example/template row.
$ws.Range("D2").Value = "{.name} is {.number} years old this year"
$ws.Range("E2").Value = "\{.name\} ignored，{.name}"

# Move the sheet's active selection to E2 (was E6).
$ws.Range("E2").Select()
